# Daily attendance processing - reorder "Recorded By" (column G) entries.
# Rule observed from the source data:
#  - If the value starts with "System", move "System" to the end of the list.
#  - Else if the value ends with lowercase "system", rotate the list right by
#    one (the trailing "system" moves to the front).
#  - Otherwise leave the value unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

function Transform-Recorder($val) {
    if ($null -eq $val) { return $val }

    $parts = $val.Split(",")
    $items = @()
    foreach ($p in $parts) {
        $items += $p.Trim()
    }
    $n = $items.Length
    if ($n -le 1) { return $val }

    $first = $items[0]
    $last = $items[$n - 1]

    if ($first.Equals("System")) {
        $rest = @()
        for ($i = 1; $i -lt $n; $i++) { $rest += $items[$i] }
        $rest += "System"
        return ($rest -join ", ")
    } elseif ($last.Equals("system")) {
        $rotated = @($last)
        for ($i = 0; $i -lt ($n - 1); $i++) { $rotated += $items[$i] }
        return ($rotated -join ", ")
    } else {
        return $val
    }
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $val -ne "") {
        $newVal = Transform-Recorder $val
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
